$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the formulas (decoder sizing spreadsheet formula corrections)
$ws.Range("F17").Formula = "=((M7-H7)*`$B`$21+H7*`$C`$21)*LN(2)"
$ws.Range("F18").Formula = '=(COUNTIF(C8:M8,"Nand3")*$H$2+COUNTIF(C8:M8,"Nand2")*$H$3+COUNTIF(C8:M8,"Inv")*$H$4)*LN(2)'
$ws.Range("B21").Formula = "=POWER(B20,(1/M7))"

# Update the active cell selection to match the saved view state
$ws.Range("F22").Select()

$wb.Application.Calculate()
